$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Skill / Frequency values for rows 2-70 (row, Skill, Frequency)
$skillData = @(
    @(2, '3Ds Max', 10),
    @(3, 'Maya', 9),
    @(4, 'Photoshop', 6),
    @(5, 'Cinema 4D', 5),
    @(6, 'Design', 4),
    @(7, '3D Design', 3),
    @(8, '3D Compositing', 3),
    @(9, 'Simulation Artist', 3),
    @(10, 'Houdini', 3),
    @(11, 'FX Artist', 3),
    @(12, 'Fume FX', 3),
    @(13, 'Time Management', 3),
    @(14, 'Technical Skills', 3),
    @(15, 'Particle Simulation', 3),
    @(16, '3D Generalist', 3),
    @(17, 'texturing', 3),
    @(18, '3D Modeler', 3),
    @(19, 'Mudbox', 3),
    @(20, 'Rigging', 3),
    @(21, '3D max', 2),
    @(22, 'CorelDraw', 2),
    @(23, 'Visual Designer', 2),
    @(24, '2D Design', 2),
    @(25, 'Visualiser', 2),
    @(26, 'Sketchup', 2),
    @(27, 'Illustrator', 2),
    @(28, 'AutoCAD', 2),
    @(29, 'Texturing', 2),
    @(30, 'Interiors', 2),
    @(31, 'Sewing', 2),
    @(32, 'CAD', 2),
    @(33, 'Visual Effects', 2),
    @(34, 'VRAY', 2),
    @(35, '3D Character Animation', 2),
    @(36, 'Art', 2),
    @(37, 'Unity3D', 2),
    @(38, '3D Maya', 2),
    @(39, 'Lumion', 1),
    @(40, 'Google Sketchup', 1),
    @(41, 'REVIT', 1),
    @(42, 'Execution', 1),
    @(43, 'Project Management', 1),
    @(44, 'Rendering', 1),
    @(45, 'Lighting', 1),
    @(46, 'New Product', 1),
    @(47, '3D Modeling', 1),
    @(48, 'Project Coordination', 1),
    @(49, 'Lead Generation', 1),
    @(50, 'Project Sales', 1),
    @(51, 'Bdm', 1),
    @(52, 'Business Development Management', 1),
    @(53, 'VAVE', 1),
    @(54, '3D Cad', 1),
    @(55, 'UG NX', 1),
    @(56, 'Teamcenter', 1),
    @(57, 'Graphics', 1),
    @(58, '3D', 1),
    @(59, 'Kaizen', 1),
    @(60, 'Time management', 1),
    @(61, 'Venture capital', 1),
    @(62, 'Architecture', 1),
    @(63, 'Adobe Premiere Pro', 1),
    @(64, 'VFX', 1),
    @(65, 'Adobe After Effects', 1),
    @(66, 'Compositing', 1),
    @(67, '3D Graphics', 1),
    @(68, 'UX', 1),
    @(69, '3D Animation', 1),
    @(70, 'Heavy Engineering', 1)
)

foreach ($entry in $skillData) {
    $row = $entry[0]
    $skill = $entry[1]
    $freq = $entry[2]
    $ws.Cells.Item($row, 2).Value = $skill
    $ws.Cells.Item($row, 3).Value = $freq
}

# Remove the now-unused trailing rows 71-80 so the sheet ends at row 70
$ws.Range("A71:C80").EntireRow.Delete() | Out-Null
